# -----------------------------------------------------------------------
# Applies the documented edit to R_templates/testreport_template.docx:
#  1. Split the "R Markdown" heading run in two and drop a fresh _GoBack
#     bookmark between "R " and "Markdown" (mirrors Word's own behaviour
#     of moving the _GoBack mark to the last edit point).
#  2. Remove the old trailing _GoBack bookmark at the end of the doc.
#  3. Heading 2 style: let "space after" fall back to the inherited 200
#     twips value (drop the explicit after="0").
#  4. Heading 5 style: add a page break before, tighten/convert the
#     spacing to an exact 1pt line, and shrink+recolor the run to a
#     near-invisible white 8pt marker (used as a page-break spacer).
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1 & 2: bookmarks -----------------------------------------------------

# Remove the stale _GoBack bookmark Word left at the end of the document.
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

# Re-create _GoBack at the point the author was last editing: between
# "R " and "Markdown" in the first heading. Word splits the run and drops
# the bookmark in automatically.
$headingRange = $d.Content
$headingRange.Find.Execute("R Markdown", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($headingRange.Start + 2, $headingRange.Start + 2)
$d.Bookmarks.Add("_GoBack", $splitPoint)

# --- 3: Heading 2 spacing ---------------------------------------------------

$heading2 = $d.Styles("Heading2")
$heading2.ParagraphFormat.SpaceAfter = 10

# --- 4: Heading 5 style ------------------------------------------------------

$heading5 = $d.Styles("Heading5")
$heading5.ParagraphFormat.PageBreakBefore = $true
$heading5.ParagraphFormat.SpaceBefore = 12
$heading5.ParagraphFormat.SpaceAfter = 0
$heading5.ParagraphFormat.LineSpacingRule = 4
$heading5.ParagraphFormat.LineSpacing = 1
$heading5.Font.Size = 8
$heading5.Font.TextColor.ObjectThemeColor = 12

Write-Output "edit complete"
